$d = $word.ActiveDocument

$replacements = @(
    @("979÷5=195, 4", "806÷8=100, 6"),
    @("256÷6=42, 4", "114÷5=22, 4"),
    @("470÷4=117, 2", "898÷4=224, 2"),
    @("473÷4=118, 1", "892÷2=446, 0"),
    @("301÷7=43, 0", "929÷2=464, 1"),
    @("795÷3=265, 0", "739÷4=184, 3"),
    @("802÷4=200, 2", "939÷8=117, 3"),
    @("810÷9=90, 0", "883÷8=110, 3"),
    @("286÷7=40, 6", "346÷2=173, 0"),
    @("715÷8=89, 3", "246÷5=49, 1"),
    @("512÷3=170, 2", "925÷8=115, 5"),
    @("345÷9=38, 3", "584÷2=292, 0"),
    @("516÷2=258, 0", "641÷5=128, 1"),
    @("218÷6=36, 2", "110÷8=13, 6"),
    @("164÷7=23, 3", "924÷7=132, 0"),
    @("226÷7=32, 2", "621÷7=88, 5"),
    @("449÷6=74, 5", "196÷9=21, 7"),
    @("314÷8=39, 2", "577÷9=64, 1"),
    @("825÷5=165, 0", "496÷7=70, 6"),
    @("236÷7=33, 5", "834÷4=208, 2"),
    @("284÷9=31, 5", "878÷5=175, 3"),
    @("509÷9=56, 5", "558÷3=186, 0"),
    @("967÷3=322, 1", "490÷8=61, 2"),
    @("627÷5=125, 2", "855÷6=142, 3"),
    @("345÷4=86, 1", "956÷5=191, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
